$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the search version value in C2 (1PSEARCHV3 -> 1PSEARCHV4)
$ws.Range("C2").Value = "1PSEARCHV4"

# Update the selected/active cell to C2
$ws.Range("C2").Select()
